$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2962962962962963
$ws.Range("C2").Value = 0.3333333333333333
$ws.Range("P2").Value = 0.2592592592592592
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.5
$ws.Range("D6").Value = 0.05263157894736842
$ws.Range("F6").Value = 0.05263157894736842
$ws.Range("J6").Value = 0.4210526315789473
$ws.Range("Q6").Value = 0.1578947368421053
$ws.Range("S6").Value = 0.3157894736842105
$ws.Range("B7").Value = 0.05
$ws.Range("F7").Value = 0.05
$ws.Range("J7").Value = 0.05
$ws.Range("O7").Value = 0.05
$ws.Range("Q7").Value = 0.1
$ws.Range("R7").Value = 0.05
$ws.Range("S7").Value = 0.65
$ws.Range("B8").Value = 0.131578947368421
$ws.Range("D8").Value = 0.02631578947368421
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.07894736842105263
$ws.Range("Q8").Value = 0.05263157894736842
$ws.Range("R8").Value = 0.02631578947368421
$ws.Range("S8").Value = 0.631578947368421
$ws.Range("B9").Value = 0.05882352941176471
$ws.Range("F9").Value = 0.05882352941176471
$ws.Range("O9").Value = 0.05882352941176471
$ws.Range("Q9").Value = 0.2352941176470588
$ws.Range("S9").Value = 0.5882352941176471
$ws.Range("B10").Value = 0.1395348837209302
$ws.Range("F10").Value = 0.05813953488372093
$ws.Range("J10").Value = 0.09302325581395349
$ws.Range("O10").Value = 0.01162790697674419
$ws.Range("Q10").Value = 0.1162790697674419
$ws.Range("R10").Value = 0.01162790697674419
$ws.Range("S10").Value = 0.5697674418604651
$ws.Range("G11").Value = 0.2647058823529412
$ws.Range("J11").Value = 0.05882352941176471
$ws.Range("K11").Value = 0.2941176470588235
$ws.Range("L11").Value = 0.2941176470588235
$ws.Range("S11").Value = 0.08823529411764706
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.25
$ws.Range("F15").Value = 0.05555555555555555
$ws.Range("H15").Value = 0.05555555555555555
$ws.Range("I15").Value = 0.05555555555555555
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.05555555555555555
$ws.Range("O15").Value = 0.05555555555555555
$ws.Range("F16").Value = 0.08333333333333333
$ws.Range("H16").Value = 0.25
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("K16").Value = 0.1666666666666667
$ws.Range("S16").Value = 0.08333333333333333
$ws.Range("F17").Value = 0.09523809523809523
$ws.Range("H17").Value = 0.09523809523809523
$ws.Range("I17").Value = 0.2380952380952381
$ws.Range("J17").Value = 0.2857142857142857
$ws.Range("K17").Value = 0.04761904761904762
$ws.Range("O17").Value = 0.1428571428571428
$ws.Range("S17").Value = 0.09523809523809523
$ws.Range("F19").Value = 0.02127659574468085
$ws.Range("H19").Value = 0.2340425531914894
$ws.Range("I19").Value = 0.07801418439716312
$ws.Range("J19").Value = 0.3120567375886525
$ws.Range("K19").Value = 0.1418439716312057
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("O19").Value = 0.06382978723404255
$ws.Range("S19").Value = 0.1276595744680851